$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.015.10'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.779.75'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.49'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5358'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('E8').Value = '  -3.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07417'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.49'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.091'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.40'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.082'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.214'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '1.778.15'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '88.54'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.71%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001052'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06483'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.38'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.901'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('D23').Value = '28.026.80'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.09'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.089'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.97'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.27'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Value = '1.979.99'
$ws.Range('E28').Value = '  -1.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.283'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.65%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '120.02'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.094'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1040'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.655'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.497'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.2238'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.72%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06357'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02261'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.968'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.451'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6158'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.436'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.87%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.178'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.30%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.004'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.19'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.670'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5750'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.99%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '125.67'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.28%  '
$ws.Range('E49').Value = '  +4.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.925'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06825'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.00%  '
